# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-column labels (row 1) to the
#   version-specific "_FV2210" / "_FV2304" labels.
# - Turn the A1:U58 range into a real Excel table (Table1).
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row shared strings in place -----------------
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$oldNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldSuffixCols.Length; $i++) {
    $addr = $oldSuffixCols[$i] + "1"
    $ws.Range($addr).Value2 = $oldNames[$i] + "_FV2210"
}

$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $newSuffixCols.Length; $i++) {
    $addr = $newSuffixCols[$i] + "1"
    $ws.Range($addr).Value2 = $oldNames[$i] + "_FV2304"
}

# --- 2) Convert the range into a table ---------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
